# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for "Nectarín" (Vega Modelo de Temuco) right
# before the existing row 590, shifting all subsequent rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 590..593; Excel copies formatting (incl. the date
# number-format on column D) from the row above, same as a normal
# "Insert Copied Cells" / "Insert Sheet Rows" operation in the UI.
$ws.Rows("590:593").Insert()

# Shared/common column values for this block of the dataset.
$marketId    = 10
$market      = "Vega Modelo de Temuco"
$region      = "La Araucanía"
$codreg      = 9
$tipo        = "Fruta"
$productoId  = 100103
$producto    = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria   = "Nectarín"

function Set-Row($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 1).Value  = $marketId
    $ws.Cells.Item($Row, 2).Value  = $market
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $productoId
    $ws.Cells.Item($Row, 8).Value  = $producto
    $ws.Cells.Item($Row, 9).Value  = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

Set-Row 590 44918 "Artic Star" "Primera" 150 20000 20000 20000 "$/bandeja 18 kilos granel" "Región de O'Higgins" 1111 18
Set-Row 591 44918 "Artic Star" "Primera" 3 450000 450000 450000 "$/bins (420 kilos)" "Región de O'Higgins" 1071 420
Set-Row 592 44918 "Artic Star" "Segunda" 170 14000 15000 14529 "$/bandeja 18 kilos granel" "Región de O'Higgins" 807 18
Set-Row 593 44918 "Early Glo" "Primera" 3 500000 500000 500000 "$/bins (420 kilos)" "Región de O'Higgins" 1190 420
